$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'98.897.97"
$ws.Range("E2").Value = "  +1.49%  "

$ws.Range("D3").Value = "'3.392.74"
$ws.Range("E3").Value = "  +8.71%  "

$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").Value = "'261.37"
$ws.Range("E5").Value = "  +8.67%  "

$ws.Range("D6").Value = "'634.17"
$ws.Range("E6").Value = "  +3.76%  "

$ws.Range("D7").Value = "'1.40"
$ws.Range("E7").Value = "  +25.78%  "

$ws.Range("E8").Value = "  +3.11%  "

$ws.Range("E9").Value = "  -0.10%  "

$ws.Range("D10").Value = "'0.885"
$ws.Range("E10").Value = "  +12.17%  "

$ws.Range("D11").Value = "'3.391.03"
$ws.Range("E11").Value = "  +8.70%  "

$ws.Range("E12").Value = "  +1.56%  "

$ws.Range("D13").Value = "'98.749.48"
$ws.Range("E13").Value = "  +1.99%  "

$ws.Range("D14").Value = "'36.46"
$ws.Range("E14").Value = "  +7.23%  "

$ws.Range("E15").Value = "  +3.44%  "

$ws.Range("D16").Value = "'4.005.39"
$ws.Range("E16").Value = "  +8.07%  "

$ws.Range("E17").Value = "  +3.45%  "

$ws.Range("D18").Value = "'3.380.04"
$ws.Range("E18").Value = "  +8.23%  "

$ws.Range("D19").Value = "'3.61"
$ws.Range("E19").Value = "  +0.32%  "

$ws.Range("D20").Value = "'15.33"
$ws.Range("E20").Value = "  +4.89%  "

$ws.Range("D21").Value = "'497.64"
$ws.Range("E21").Value = "  -2.98%  "

$ws.Range("E22").Value = "  +9.48%  "

$ws.Range("E23").Value = "  +9.58%  "

$ws.Range("D24").Value = "'9.44"
$ws.Range("E24").Value = "  +6.47%  "

$ws.Range("D25").Value = "'5.76"
$ws.Range("E25").Value = "  +4.03%  "

$ws.Range("D26").Value = "'89.41"
$ws.Range("E26").Value = "  +3.33%  "

$ws.Range("D27").Value = "'12.17"
$ws.Range("E27").Value = "  +3.96%  "

$ws.Range("D29").Value = "'0.284"
$ws.Range("E29").Value = "  +19.95%  "

$ws.Range("D30").Value = "'0.204"
$ws.Range("E30").Value = "  +16.67%  "

$ws.Range("E31").Value = "  -0.36%  "

$ws.Range("E32").Value = "  +6.48%  "

$ws.Range("D33").Value = "'0.999"
$ws.Range("E33").Value = "  +18.12%  "

$ws.Range("D34").Value = "'9.65"
$ws.Range("E34").Value = "  +6.79%  "

$ws.Range("D35").Value = "'28.09"
$ws.Range("E35").Value = "  +5.36%  "

$ws.Range("E36").Value = "  +0.47%  "

$ws.Range("E37").Value = "  -0.05%  "

$ws.Range("D38").Value = "'2.00"
$ws.Range("E38").Value = "  +6.58%  "

$ws.Range("E39").Value = "  +8.28%  "

$ws.Range("D40").Value = "'505.55"
$ws.Range("E40").Value = "  +2.55%  "

$ws.Range("E41").Value = "  +2.54%  "

$ws.Range("D42").Value = "'1.29"
$ws.Range("E42").Value = "  +2.96%  "

$ws.Range("D43").Value = "'3.74"
$ws.Range("E43").Value = "  +4.06%  "

$ws.Range("E44").Value = "  +5.86%  "

$ws.Range("D45").Value = "'0.794"
$ws.Range("E45").Value = "  +14.26%  "

$ws.Range("E46").Value = "  -0.02%  "

$ws.Range("D47").Value = "'160.63"
$ws.Range("E47").Value = "  -1.31%  "

$ws.Range("D48").Value = "'1.97"
$ws.Range("E48").Value = "  +2.69%  "

$ws.Range("D49").Value = "'0.847"
$ws.Range("E49").Value = "  +15.72%  "

$ws.Range("D50").Value = "'4.72"
$ws.Range("E50").Value = "  +7.58%  "

$ws.Range("D51").Value = "'46.63"
$ws.Range("E51").Value = "  +4.75%  "
